$d = $word.ActiveDocument

# --- 1. "Cursor trail" -> "Responsiveness" ---
$d.Content.Find.Execute("Cursor trail", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Responsiveness", 2) | Out-Null

# --- 2. "Add flip side" -> "Add photo picture, random" (will absorb the old "Parallax" para) ---
$d.Content.Find.Execute("Add flip side", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Add photo picture, random", 2) | Out-Null

# --- 3. Remove the (now redundant) "Parallax" paragraph ---
$d.Content.Find.Execute("Parallax", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 2) | Out-Null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "") {
        $p.Range.Delete() | Out-Null
        break
    }
}

# --- 4. "Give credit" -> "Stop from animating everytime" (the paragraph 5 occurrence) ---
$d.Content.Find.Execute("Give credit", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Stop from animating everytime", 2) | Out-Null

# --- 5. Remove the _GoBack bookmark currently on "Flatten images" ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 6. Re-insert "Give credit" as a new paragraph right after "Flatten images" ---
$flatten = $d.Paragraphs.Item($d.Paragraphs.Count)
$flatten.Range.InsertParagraphAfter()
$newGiveCredit = $d.Paragraphs.Item($d.Paragraphs.Count)
$newGiveCredit.Range.Text = "Give credit"

# --- 7. Insert the new "Top panel " paragraph right after "Responsiveness" ---
$resp = $d.Paragraphs.Item(1)
$resp.Range.InsertParagraphAfter()
$topPanel = $d.Paragraphs.Item(2)
$topPanel.Range.Text = "Top panel "

# --- 8. Re-add the _GoBack bookmark at the end of the "Responsiveness" paragraph ---
$resp = $d.Paragraphs.Item(1)
$bmRange = $resp.Range.Duplicate
$bmRange.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- 9. Append two empty paragraphs at the very end of the document ---
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.Text = [char]13

$endRng2 = $d.Content
$endRng2.Collapse(0)
$endRng2.Text = [char]13

Write-Output "done"
